$d = $word.ActiveDocument

# The diff turns three single-<w:t> runs (English translations of the
# "Objetivos"/"Programa resumido"/"Programa" sections) into runs that are
# split into several <w:t> pieces separated by manual line breaks (<w:br/>).
# Word's "Find & Replace" treats "^l" in the replacement string as a manual
# line break, so we use Find.Execute with wildcards off and Replace=2
# (wdReplaceAll) to insert each break at the right spot.

# Paragraph: "General- Show the Analytical Chemistry ... sensible.Specifics: - By completing ..."
# -> "General" / "- Show the Analytical Chemistry ... sensible." / "Specifics: " / "- By completing ..."
$d.Content.Find.Execute(
    "General- Show the Analytical Chemistry",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "General^l- Show the Analytical Chemistry",
    2)

$d.Content.Find.Execute(
    "sensible.Specifics: - By completing",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "sensible.^lSpecifics: ^l- By completing",
    2)

# Paragraph: "- Theoretical bases of analytical chemistry ... qualitative analysis.- Fundamentals of analysis titrimetry ..."
# -> "- Theoretical bases ... qualitative analysis." / "- Fundamentals of analysis titrimetry ..."
$d.Content.Find.Execute(
    "qualitative analysis.- Fundamentals of analysis titrimetry",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "qualitative analysis.^l- Fundamentals of analysis titrimetry",
    2)

# Paragraph: "- Theoretical Bases of Qualitative Analysis ... by Vogel.- Foundations of Quantitative Analysis ..."
# -> "- Theoretical Bases of Qualitative Analysis ... by Vogel." / "- Foundations of Quantitative Analysis ..."
$d.Content.Find.Execute(
    "by Vogel.- Foundations of Quantitative Analysis",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "by Vogel.^l- Foundations of Quantitative Analysis",
    2)
